# Add newly uploaded venue/county rows to the existing list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("DEFY Pairc Mhuire, Ardee", "Louth"),
    @("Athleague", "Roscommon"),
    @("FBD Semple Stadium", "Tipperary"),
    @("Manguard Park", "Kildare"),
    @("Castleblayney", "Monaghan"),
    @("Crossmaglen", "Armagh"),
    @("Ballygunner", "Waterford"),
    @("Bansha", "Tipperary"),
    @("Grant Heating St Brendan's Park, Birr", "Offaly"),
    @("Middletown", "Armagh"),
    @("Kent Park, Ballydoogan, Sligo", "Sligo"),
    @("Fethard Town Park, Grass Pitch", "Tipperary"),
    @("Lavey", "Derry"),
    @("SETU Arena, Carriganore", "Waterford"),
    @("St Joseph's Glenavy, Antrim", "Antrim"),
    @("Páirc Éanna, Glengormley", "Antrim"),
    @("O Neill Park, Dungannon", "Tyrone"),
    @("Shane McGettigan Park, Drumshanbo", "Leitrim"),
    @("Inniskeen", "Monaghan"),
    @("Páirc Uí Rinn", "Cork"),
    @("Healy Park, Omagh", "Tyrone")
)

$startRow = 73
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# Autofit the two data columns, matching the bestFit column widths Excel
# records after the paste (column A ~52.14, column B ~13.14 chars).
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null

# Restore the selection Excel left on the sheet after the edit.
$ws.Range("C14").Select() | Out-Null
